$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# Title (appears twice: main heading and bolded text near the end)
Replace-Text "Play Lost Riches of El Dorado Free Slot Game Review" "Play Lost Riches of El Dorado for Free"

# "What we like" bullet list
Replace-Text "Engaging adventure theme centered around El Dorado" "Engaging adventure theme"
Replace-Text "Excellent quality graphics and atmosphere" "Excellent graphics"
Replace-Text "Rich in exciting features like free spins, cash bonus, and gamble feature" "Wide range of betting options"
Replace-Text "Suitable betting range for beginners" "Exciting bonus features"

# "What we don't like" bullet list
Replace-Text "Purchasing free spins option is quite expensive" "Purchasing free spins can be expensive"
Replace-Text "No progressive jackpot" "Limited number of paylines"

# Final italic summary line
Replace-Text "Read our unbiased review of Lost Riches of El Dorado slot game and play for free. Engaging adventure theme, excellent graphics, and rich in exciting features." "Read our review of Lost Riches of El Dorado and play this exciting adventure-themed slot game for free."
